$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column right after "ESTADO" (column A) for the new "ID"
# header used by the code-scanning feature, shifting the existing
# headers (NOMBRES Y APELLIDOS, CEDULA, ... OBSERVACIONES) one column
# to the right.
$ws.Columns("B:B").Insert()

# Give the new header the same look (bold / centered / bordered) as
# the rest of the header row by copying the formatting from the
# neighboring header cell, then set its text.
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B1").Value = "ID"
